$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns J, K, L (they are no longer part of the data range)
$ws.Range("J1:L12").Delete()

# Row 1 header values stay 0..7 in B1:I1 (unchanged, nothing to do)

# Update data rows 2-12, columns G, H, I (values shifted from old H/J, and I unchanged
# except row 11 which also got a new I value) plus row 11's B:F direct value edits.

$ws.Range("G2").Value = 0.127
$ws.Range("H2").Value = 0.219
$ws.Range("I2").Value = 0.295

$ws.Range("G3").Value = 0.079
$ws.Range("H3").Value = 0.181
$ws.Range("I3").Value = 0.356

$ws.Range("G4").Value = 0.105
$ws.Range("H4").Value = 0.148
$ws.Range("I4").Value = 0.437

$ws.Range("G5").Value = 0.152
$ws.Range("H5").Value = 0.133
$ws.Range("I5").Value = 0.511

$ws.Range("G6").Value = 0.127
$ws.Range("H6").Value = 0.108
$ws.Range("I6").Value = 0.48

$ws.Range("G7").Value = 0.082
$ws.Range("H7").Value = 0.078
$ws.Range("I7").Value = 0.332

$ws.Range("G8").Value = 0.054
$ws.Range("H8").Value = 0.059
$ws.Range("I8").Value = 0.283

$ws.Range("G9").Value = 0.023
$ws.Range("H9").Value = 0.043
$ws.Range("I9").Value = 0.28

$ws.Range("G10").Value = 0.024
$ws.Range("H10").Value = 0.042
$ws.Range("I10").Value = 0.279

$ws.Range("B11").Value = 1.002
$ws.Range("C11").Value = 0.98
$ws.Range("D11").Value = -2.975
$ws.Range("E11").Value = 0.992
$ws.Range("F11").Value = 0.081
$ws.Range("G11").Value = 0.025
$ws.Range("H11").Value = 0.042
$ws.Range("I11").Value = 0.282

$ws.Range("G12").Value = 0.023
$ws.Range("H12").Value = 0.041
$ws.Range("I12").Value = 0.28
